# Apply variety-tracking updates on the "2€" sheet of the commemorative
# Germany coin-collection workbook: several K:O "varieties owned" flags
# move from 0 to 1 (newly acquired varieties), as reflected by the
# updated P-column "Can exchange" formulas already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> list of columns (K..O) that flip from 0 to 1.
$changes = @{
    3  = @("L", "M", "N", "O")
    4  = @("N", "O")
    5  = @("L", "M", "N", "O")
    7  = @("M", "N", "O")
    8  = @("M")
    9  = @("K", "M")
    10 = @("L", "M")
    11 = @("L", "M")
    12 = @("L", "O")
    15 = @("K", "M", "N", "O")
    16 = @("L", "M", "N", "O")
    17 = @("K", "L", "M", "N")
    20 = @("L", "M", "N", "O")
    21 = @("M")
}

foreach ($row in $changes.Keys) {
    foreach ($col in $changes[$row]) {
        $ws.Range("$col$row").Value = 1
    }
}

# Reflect the last on-screen selection recorded for the sheet (the
# bottom-right frozen pane) after making the edits.
$ws.Activate()
$ws.Range("P27").Select()
